$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText {
    param($table, $row, $col, $newText)
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    # Exclude the trailing cell-mark character so only the visible text
    # (and its paragraph mark) is targeted; this keeps the existing run
    # formatting (font/size) intact while swapping the text content.
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $textRange.Text = $newText
}

# Row 1 (first data row)
Set-CellText $t 1 1 "256×2=512"
Set-CellText $t 1 2 "538×9=4842"
Set-CellText $t 1 3 "867×6=5202"
Set-CellText $t 1 4 "316×2=632"
Set-CellText $t 1 5 "215×6=1290"

# Row 5
Set-CellText $t 5 1 "572×8=4576"
Set-CellText $t 5 2 "607×6=3642"
Set-CellText $t 5 3 "737×5=3685"
Set-CellText $t 5 4 "564×3=1692"
Set-CellText $t 5 5 "877×5=4385"

# Row 10
Set-CellText $t 10 1 "443×6=2658"
Set-CellText $t 10 2 "911×6=5466"
Set-CellText $t 10 3 "847×2=1694"
Set-CellText $t 10 4 "469×7=3283"
Set-CellText $t 10 5 "953×7=6671"

# Row 15
Set-CellText $t 15 1 "356×7=2492"
Set-CellText $t 15 2 "622×8=4976"
Set-CellText $t 15 3 "189×9=1701"
Set-CellText $t 15 4 "619×2=1238"
Set-CellText $t 15 5 "909×8=7272"

# Row 20
Set-CellText $t 20 1 "745×8=5960"
Set-CellText $t 20 2 "143×5=715"
Set-CellText $t 20 3 "784×9=7056"
Set-CellText $t 20 4 "565×9=5085"
Set-CellText $t 20 5 "838×8=6704"

Write-Host "Done updating table cells."
